$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = 45955
$ws.Range("B2").Value = 97.59
$ws.Range("C2").Value = 93.08
$ws.Range("D2").Value = 87.70999999999999
$ws.Range("E2").Value = 84.47
$ws.Range("F2").Value = 95.40000000000001
$ws.Range("G2").Value = 96.39
$ws.Range("H2").Value = 99.51000000000001
$ws.Range("I2").Value = 108.14
$ws.Range("J2").Value = 117.41
$ws.Range("K2").Value = 104.38
$ws.Range("L2").Value = 84.5
$ws.Range("M2").Value = 56.79
$ws.Range("N2").Value = 35.71
$ws.Range("O2").Value = 28.96
$ws.Range("P2").Value = 25.02
$ws.Range("Q2").Value = 25.02
$ws.Range("R2").Value = 25.02
$ws.Range("S2").Value = 47.04
$ws.Range("T2").Value = 88.59
$ws.Range("U2").Value = 113.12
$ws.Range("V2").Value = 126.49
$ws.Range("W2").Value = 114.04
$ws.Range("X2").Value = 106.71
$ws.Range("Y2").Value = 99.78
$ws.Range("Z2").Value = 81.7
$ws.Range("AB2").Value = 111.76
$ws.Range("AD2").Value = 120.26
$ws.Range("AF2").Value = 110.9
$ws.Range("AG2").Value = "11h-17h"
